$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing rows down
$ws.Rows.Item(2).Insert()

# Copy previous "row 2" content (now row 3) into new row 2, with updated date
$ws.Cells.Item(2, 1).Value = "07-11-2025"
$ws.Cells.Item(2, 2).Value = $ws.Cells.Item(3, 2).Value
$ws.Cells.Item(2, 3).Value = $ws.Cells.Item(3, 3).Value
$ws.Cells.Item(2, 4).Value = $ws.Cells.Item(3, 4).Value
$ws.Cells.Item(2, 5).Value = $ws.Cells.Item(3, 5).Value
$ws.Cells.Item(2, 6).Value = $ws.Cells.Item(3, 6).Value
